# Zero Trust Research -> Zero Trust Architecture + article summary + bibliography
$d = $word.ActiveDocument

# --- 1. Split the heading run "Zero Trust Research" into "Zero Trust" + " Architecture" ---
$headingPara = $d.Paragraphs(1)
$headingRange = $headingPara.Range
$headingTextLen = $headingRange.End - $headingRange.Start
# Exclude the trailing paragraph mark from the replaced range.
$headingTextRange = $d.Range($headingRange.Start, $headingRange.End - 1)
if ($headingTextRange.Text -eq "Zero Trust Research") {
    $headingXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:r><w:t>Zero Trust</w:t></w:r><w:r><w:t xml:space="preserve"> Architecture</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $headingTextRange.InsertXML($headingXml)
}

# --- 2. Append the new content (citation, summary paragraphs, bibliography) ---
# at the very end of the main document body, right before the sectPr.
$endPos = $d.Content.End
$insertRange = $d.Range($endPos, $endPos)
$bodyXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p><w:sdt><w:sdtPr><w:id w:val="-1418405477"/><w:citation/></w:sdtPr><w:sdtContent><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> CITATION Dhi24 \l 6153 </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>(Dhiman, et al., 2024)</w:t></w:r><w:r><w:fldChar w:fldCharType="end"/></w:r></w:sdtContent></w:sdt></w:p><w:p/><w:p><w:r><w:t>This is a short summary and review of an article on a Zero Trust architecture.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">A zero trust architecture is an emerging software security architecture paradigm. A zero trust architecture can be presumed as </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>multi faceted</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>. There is no solitary technology or architecture that completely implements a zero trust</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">model. In realising this, when designing an architecture, environment specific implementation strategies must be devised. The article referenced in this document discusses such implementation strategies and their </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>locical</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> components.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t xml:space="preserve">The paper delves into a comparative analysis of zero trust systems, where </w:t></w:r><w:r><w:t>various technologies are assessed for their suitability. It discusses important parameters surrounding the importance of operational requirements over efficiency, delving into how open source software and microservices play a key part in enhancing security and rapid deployment of software, and maintenance simplification.</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>The paper further expands on the need to emphasize the inclusion of zero trust policies in future technologies and architecture types, including 5G/6G networking, edge computing deployments and further discusses intelligent zero trust applications as a security mechanism for untrusted networking components. It expands upon the use of artificial intelligence as a method of enhancing security measures in an architecture, which instigates a shift from reactive network security applications to proactive, where early detection is key.</w:t></w:r></w:p><w:p/><w:p/><w:p><w:r><w:t xml:space="preserve">The study detailed in this paper emphasises the importance of correct application of </w:t></w:r><w:r><w:t>authentication and access control approaches, where organisations constantly re-evaluate their trust in active connection points in their architecture. Since each organisation has their own unique implementation the application of a zero trust architecture will differ greatly from one use case to the next. The article elaborates on how it’s important to make use of proper encryption techniques, and segment software into smaller components (micro services). Zero trust architecture is quite an emerging technology and will evolve with further studies such as this in the near and distant future.</w:t></w:r></w:p><w:p/><w:sdt><w:sdtPr><w:id w:val="675231696"/><w:docPartObj><w:docPartGallery w:val="Bibliographies"/><w:docPartUnique/></w:docPartObj></w:sdtPr><w:sdtEndPr><w:rPr><w:rFonts w:asciiTheme="minorHAnsi" w:eastAsiaTheme="minorHAnsi" w:hAnsiTheme="minorHAnsi" w:cstheme="minorBidi"/><w:color w:val="auto"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:sdtEndPr><w:sdtContent><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Bibliography</w:t></w:r></w:p><w:sdt><w:sdtPr><w:id w:val="111145805"/><w:bibliography/></w:sdtPr><w:sdtContent><w:p><w:pPr><w:pStyle w:val="Bibliography"/><w:ind w:left="720" w:hanging="720"/><w:rPr><w:noProof/><w:kern w:val="0"/><w14:ligatures w14:val="none"/></w:rPr></w:pPr><w:r><w:fldChar w:fldCharType="begin"/></w:r><w:r><w:instrText xml:space="preserve"> BIBLIOGRAPHY </w:instrText></w:r><w:r><w:fldChar w:fldCharType="separate"/></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t xml:space="preserve">Dhiman, P., Saini, N., Gulzar, Y., Turaev, S., Kaur, A., Nisa, K., &amp; Hamid, Y. (2024, February 19). A Review and Comparative Analysis of Relevant Approaches of Zero Trust Network Model. </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:noProof/></w:rPr><w:t>Sensors, 24</w:t></w:r><w:r><w:rPr><w:noProof/></w:rPr><w:t>(4), 1-19.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:b/><w:bCs/><w:noProof/></w:rPr><w:fldChar w:fldCharType="end"/></w:r></w:p></w:sdtContent></w:sdt></w:sdtContent></w:sdt><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertRange.InsertXML($bodyXml)

# --- 3. Register the "Bibliography" paragraph style used above ---
$biblioStyle = $null
try {
    $biblioStyle = $d.Styles("Bibliography")
} catch {
    $biblioStyle = $null
}
if (-not $biblioStyle) {
    $biblioStyle = $d.Styles.Add("Bibliography", 1)
}
$biblioStyle.BaseStyle = "Normal"
$biblioStyle.NextParagraphStyle = "Normal"
$biblioStyle.Priority = 37
$biblioStyle.UnhideWhenUsed = $true

Write-Host "Done. Final text length: $($d.Content.End)"
